$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# --- Fix the role on rows 17-18 (Business-Process Analyst -> System Analyst) ---
$ws.Range("B17").Value = "System Analyst "
$ws.Range("B18").Value = "System Analyst "

# --- Fill in the newly logged time entries for rows 21-27 ---
$ws.Range("A21").Value = "Udarbejdelse af iterationsplan 2"
$ws.Range("B21").Value = "Project Manager"
$ws.Range("C21").Value = 43889
$ws.Range("D21").Value = 0.35416666666666669
$ws.Range("E21").Value = 0.375

$ws.Range("A22").Value = "Krydstjek af UC05 med Benjamin"
$ws.Range("B22").Value = "Business-Process Analyst"
$ws.Range("C22").Value = 43889
$ws.Range("D22").Value = 0.375
$ws.Range("E22").Value = 0.39583333333333331

$ws.Range("A23").Value = "Udarbejdelse af iterationsplan 2"
$ws.Range("B23").Value = "Project Manager"
$ws.Range("C23").Value = 43889
$ws.Range("D23").Value = 0.39583333333333331
$ws.Range("E23").Value = 0.41666666666666669

$ws.Range("A24").Value = "Review af ATD med Nicky"
$ws.Range("B24").Value = "Reviewer"
$ws.Range("C24").Value = 43889
$ws.Range("D24").Value = 0.42708333333333331
$ws.Range("E24").Value = 0.45833333333333331

$ws.Range("A25").Value = "Lavet opgaver til Whiteboard"
$ws.Range("B25").Value = "Project Manager"
$ws.Range("C25").Value = 43889
$ws.Range("D25").Value = 0.45833333333333331
$ws.Range("E25").Value = 0.47916666666666669

$ws.Range("A26").Value = "Rettet AD06"
$ws.Range("B26").Value = "Business-Process Analyst"
$ws.Range("C26").Value = 43889
$ws.Range("D26").Value = 0.51388888888888895
$ws.Range("E26").Value = 0.54166666666666663

$ws.Range("A27").Value = "Krydstjek af AD05 med Marc"
$ws.Range("B27").Value = "Business-Process Analyst"
$ws.Range("C27").Value = 43889
$ws.Range("D27").Value = 0.54166666666666663
$ws.Range("E27").Value = 0.58333333333333337

# --- Re-apply the role/date/time data validation so it also covers the new rows ---
$roleOld = $ws.Range("B22:B1048576")
$roleOld.Validation.Delete()
$roleOld.Validation.Add(3, 1, 1, "=Roller")
$roleNew = $ws.Range("B3:B21")
$roleNew.Validation.Delete()
$roleNew.Validation.Add(3, 1, 1, "=Roller")

$dateOld = $ws.Range("C22:C1048576")
$dateOld.Validation.Delete()
$dateOld.Validation.Add(4, 1, 1, "43881", "43908")
$dateNew = $ws.Range("C1:C21")
$dateNew.Validation.Delete()
$dateNew.Validation.Add(4, 1, 1, "43881", "43908")

$timeOld = $ws.Range("D22:D1048576")
$timeOld.Validation.Delete()
$timeOld.Validation.Add(5, 1, 1, "0.333333333333333", "0.708333333333333")
$timeNew = $ws.Range("D1:D21")
$timeNew.Validation.Delete()
$timeNew.Validation.Add(5, 1, 1, "0.333333333333333", "0.708333333333333")

# --- Move the active selection to A2 ---
$ws.Range("A2").Select()
